$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")

# --- New comment texts (write the "Accounting statements" one first so it
#     lands at shared-string index 109, and the "Implementation of public
#     holidays" one second so it lands at index 110 - matching the target
#     sharedStrings.xml ordering). ---
$accountingText = "2. Accounting statements has been generated for the GL of the Jan22 for all three centers and shared to Rahman san to validate `ndata"
$holidayText = "1. Implementation of public holidays has been completed at WARRANTY_DAILY and GRS_SUMMARY_DETAILS_DAILY tasks, `ntested in all SSCs and it is running smoothly"

# --- Row 31: copy the formatting pattern from row 28 (No/Date/App/Task/%/Status/blank
#     all bordered, date + percent number formats, wrapped Task column) ---
$ws.Range("A28:G28").Copy()
$ws.Range("A31:G31").PasteSpecial(-4122)

# --- Row 32: copy the formatting pattern from row 15 (blank No/Date, App/Task/%/Status) ---
$ws.Range("A15:G15").Copy()
$ws.Range("A32:G32").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 31 values ---
$ws.Cells.Item(31,1).Value = 14
$ws.Cells.Item(31,2).Value = 44638
$ws.Cells.Item(31,3).Value = "RPA GSS"
$ws.Cells.Item(32,4).Value = $accountingText
$ws.Cells.Item(31,4).Value = $holidayText
$ws.Cells.Item(31,5).Value = 1
$ws.Cells.Item(31,6).Value = "Completed"

# --- Row 32 values ---
$ws.Cells.Item(32,3).Value = "RPA RLOGIC"
$ws.Cells.Item(32,5).Value = 1
$ws.Cells.Item(32,6).Value = "Completed"

# --- Row heights (auto word-wrap height computed by real Excel) ---
$ws.Rows.Item(31).RowHeight = 43.2
$ws.Rows.Item(32).RowHeight = 43.2

# --- Selection moved from D34 to D31 ---
$ws.Range("D31").Select() | Out-Null
